$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.261112689971924
$ws.Range("B1").Value = 2.224258899688721
$ws.Range("C1").Value = 4.381772518157959
$ws.Range("D1").Value = 3.003652811050415
$ws.Range("E1").Value = 1.029847502708435
